$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2023" column (T) by copying the formatting of the existing
#     last data column (S) for rows 2-6 (row 1 has no data past column C/F) ---
$ws.Range("S2:S6").Copy() | Out-Null
$ws.Range("T2:T6").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

# --- Populate the new 2023 values ---
$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 263951
$ws.Range("T5").Value = 3.7
$ws.Range("T6").Value = 32.299999999999997

# --- Row 2 now gets an explicit custom row height ---
$ws.Rows.Item(2).RowHeight = 16.5

# --- Row 3: years 2020-2023 (Q3:T3) now share the same centred style as the
#     rest of the year header row (D3:P3) instead of their previous style ---
$ws.Range("P3").Copy() | Out-Null
$ws.Range("Q3:T3").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

# Restore / set the year values for Q3:T3 after the format paste
$ws.Range("Q3").Value = 2020
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023
